$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.850.11'
$ws.Range("E2").Value = '  +1.56%  '
$ws.Range("D3").Value = '3.701.21'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '614.42'
$ws.Range("E5").Value = '  +6.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '186.22'
$ws.Range("E6").Value = '  +4.84%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.633'
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.710'
$ws.Range("E9").Value = '  -0.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.159'
$ws.Range("E10").Value = '  -3.43%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '55.39'
$ws.Range("E11").Value = '  +5.60%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000286'
$ws.Range("E12").Value = '  -4.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.46'
$ws.Range("E13").Value = '  -1.17%  '
$ws.Range("D14").Value = '4.297.76'
$ws.Range("E14").Value = '  +0.19%  '
$ws.Range("D15").Value = '3.701.17'
$ws.Range("E15").Value = '  -0.70%  '
$ws.Range("B16").Value = 'TRON'
$ws.Range("C16").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.126'
$ws.Range("E16").Value = '  -0.39%  '
$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.17'
$ws.Range("E17").Value = '  -0.73%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.84'
$ws.Range("E18").Value = '  -0.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.12'
$ws.Range("E19").Value = '  -1.40%  '
$ws.Range("D20").Value = '68.516.42'
$ws.Range("E20").Value = '  +1.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '408.40'
$ws.Range("E21").Value = '  +0.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.57'
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '88.69'
$ws.Range("E23").Value = '  +0.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.00'
$ws.Range("E24").Value = '  -2.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.70'
$ws.Range("E25").Value = '  -0.81%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.88'
$ws.Range("E26").Value = '  +1.60%  '
$ws.Range("E27").Value = '  +1.15%  '
$ws.Range("E28").Value = '  -3.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.57'
$ws.Range("E29").Value = '  +0.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.86'
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.16'
$ws.Range("E31").Value = '  -10.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.53'
$ws.Range("E32").Value = '  -1.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.121'
$ws.Range("E33").Value = '  +2.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '626.39'
$ws.Range("E34").Value = '  +4.88%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '43.92'
$ws.Range("E35").Value = '  -0.94%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '65.38'
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("B37").Value = 'TheGraph'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.409'
$ws.Range("E37").Value = '  +2.14%  '
$ws.Range("B38").Value = 'Dai'
$ws.Range("C38").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("B39").Value = 'PEPE'
$ws.Range("C39").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D39").Value = '0.0₃0799'
$ws.Range("E39").Value = '  -12.76%  '
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.139'
$ws.Range("E41").Value = '  +2.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.00'
$ws.Range("E42").Value = '  -2.20%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0439'
$ws.Range("E43").Value = '  +0.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.60'
$ws.Range("E44").Value = '  +0.70%  '
$ws.Range("E45").Value = '  +3.03%  '
$ws.Range("D46").Value = '2.855.03'
$ws.Range("E46").Value = '  +4.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.72'
$ws.Range("E47").Value = '  +0.81%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.03'
$ws.Range("E48").Value = '  -4.44%  '
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.69'
$ws.Range("E49").Value = '  -0.68%  '
$ws.Range("B50").Value = 'ApeXProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.08'
$ws.Range("E50").Value = '  -2.52%  '
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.75'
$ws.Range("E51").Value = '  -0.47%  '